# "add links to chapters 10.8 and 11.2 to lab06"
#
# The "10.3" chapter row (row 16) has no lab links at all, so it is
# removed from the sheet entirely. Once it's gone, the "10.8" chapter
# (which becomes row 16) and the "11.2" chapter (which becomes row 18)
# each get a link recorded for lab06 -- column E -- by writing a 1 there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-empty "10.3" row; everything below shifts up by one.
$ws.Rows("16").Delete()

# lab06 (column E) links for "10.8" (now row 16) and "11.2" (now row 18).
$ws.Range("E16").Value = 1
$ws.Range("E18").Value = 1

# The "completed" column N used to just SUM the lab columns, which broke
# once a chapter could have links in more than one lab column (11.2 now
# has both C and E set). Switch it to a boolean-style IF() so it keeps
# reporting 0/1 regardless of how many labs link to a chapter.
$ws.Range("N3").Formula = "=IF(SUM(B3:M3)>0,1,0)"
$ws.Range("N4:N35").Formula = "=IF(SUM(B4:M4)>0,1,0)"

# The conditional formatting range still referred to the old N3:N36
# extent; shrink it to match the now-shorter table.
$fc = $ws.Range("N3:N36").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("N3:N35"))

# Leave the selection where the author left it.
$ws.Range("E22").Select()
